$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("report")

$ws.Range("F4").Value = "Was not able to save Interest Rates, "
$ws.Range("G4").Value = "2022-09-06 16:41:58"
